$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute(
        $old, $true, $false, $false, $false, $false, $true, 1, $false,
        $new, 2)
}

# Title (appears twice: H1 heading + bold summary run near the end)
Replace-Text "Play Empire V Free: Review of Novomatic's Online Slot Game" "Play Empire V Free - Online Slot Review"

# "What we like" bullet list
Replace-Text "Engaging vampire theme and symbols" "Well-executed vampire theme with immersive design and sound effects"
Replace-Text "Medium volatility with a simple 5x3 grid and 10 paylines" "Simple and enjoyable gameplay experience with a 5x3 grid and 10 paylines"
Replace-Text "Special symbols, including Wild and Scatter" "Medium volatility offers a good balance between wins and payouts"
Replace-Text "Bonus mode with free spins rounds" "Bonus mode with free spins rounds for increased chances of winning"

# "What we don't like" bullet list
Replace-Text "Limited number of paylines compared to other online slot games" "Limited variety of bonus features"
Replace-Text "No significant innovation or unique features" "Lacks innovative gameplay mechanics"

# Italic summary/meta description run at the very end
Replace-Text "Empire V is a vampire-themed online slot game by Novomatic with 10 paylines, special symbols, and a bonus mode with free spins rounds. Play free and read our review." "Play Empire V for free and read our review of this online slot game."
